$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2158
$ws.Range("I3").Value = 7487
$ws.Range("J3").Value = 2261
$ws.Range("C4").Value = 1821
$ws.Range("J4").Value = 512
$ws.Range("J6").Value = 2869
$ws.Range("C7").Value = 28364
$ws.Range("I7").Value = 26202
$ws.Range("J7").Value = 7962

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 97
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 116
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 292

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 61
$ws.Range("J7").Value = 239
$ws.Range("J8").Value = 500
$ws.Range("J11").Value = 110
$ws.Range("J14").Value = 27
$ws.Range("J18").Value = 93
$ws.Range("J19").Value = 263
$ws.Range("J20").Value = 164
$ws.Range("J22").Value = 16
$ws.Range("J23").Value = 74
$ws.Range("J29").Value = 446
$ws.Range("J30").Value = 30
$ws.Range("J31").Value = 54
$ws.Range("J36").Value = 120
$ws.Range("J37").Value = 267
$ws.Range("J42").Value = 304
$ws.Range("J44").Value = 67
$ws.Range("J48").Value = 73
$ws.Range("J49").Value = 47
$ws.Range("J51").Value = 106
$ws.Range("J52").Value = 195
$ws.Range("C53").Value = 362
$ws.Range("I53").Value = 296
$ws.Range("J53").Value = 76
$ws.Range("J60").Value = 54
$ws.Range("I63").Value = 205
$ws.Range("J63").Value = 35
$ws.Range("J65").Value = 207
$ws.Range("J67").Value = 292
$ws.Range("J73").Value = 71
$ws.Range("J76").Value = 116
$ws.Range("J79").Value = 239
$ws.Range("J83").Value = 190
$ws.Range("J85").Value = 377
$ws.Range("J86").Value = 45
$ws.Range("J88").Value = 82
$ws.Range("J89").Value = 83
$ws.Range("J94").Value = 61
$ws.Range("J96").Value = 89
$ws.Range("J97").Value = 53
$ws.Range("I100").Value = 39
$ws.Range("C101").Value = 28364
$ws.Range("I101").Value = 26202
$ws.Range("J101").Value = 7962

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 190

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 132
$ws.Range("J3").Value = 154
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 446

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 65
$ws.Range("J3").Value = 72
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 23
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 14
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 90
$ws.Range("J3").Value = 147
$ws.Range("J4").Value = 24
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 377

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 70
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 93

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 41
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I3").Value = 6
$ws.Range("I6").Value = 39

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 10
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 110

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 162
$ws.Range("J7").Value = 500

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 21
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 30
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 54

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 13
$ws.Range("C4").Value = 26
$ws.Range("I4").Value = 21
$ws.Range("C7").Value = 362
$ws.Range("I7").Value = 296
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 239
